$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '28.457.64'
$ws.Range("E2").Value = '  -3.72%  '
Set-TextValue $ws.Range("D3") '1.953.76'
$ws.Range("E3").Value = '  -2.48%  '
Set-TextValue $ws.Range("D4") '1.011'
$ws.Range("E4").Value = '  -0.26%  '
Set-TextValue $ws.Range("D5") '321.07'
$ws.Range("E5").Value = '  -2.54%  '
$ws.Range("E6").Value = '  -0.23%  '
Set-TextValue $ws.Range("D7") '0.4759'
$ws.Range("E7").Value = '  -4.82%  '
Set-TextValue $ws.Range("D8") '0.4066'
$ws.Range("E8").Value = '  -3.57%  '
Set-TextValue $ws.Range("D9") '53.44'
$ws.Range("E9").Value = '  -1.65%  '
Set-TextValue $ws.Range("D10") '0.08449'
$ws.Range("E10").Value = '  -6.22%  '
Set-TextValue $ws.Range("D11") '1.057'
$ws.Range("E11").Value = '  -5.38%  '
Set-TextValue $ws.Range("D12") '22.09'
$ws.Range("E12").Value = '  -5.05%  '
Set-TextValue $ws.Range("D13") '1.991.90'
$ws.Range("E13").Value = '  -2.76%  '
Set-TextValue $ws.Range("D14") '7.612'
$ws.Range("E14").Value = '  -5.20%  '
Set-TextValue $ws.Range("D15") '6.172'
$ws.Range("E15").Value = '  -4.55%  '
Set-TextValue $ws.Range("D16") '1.012'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("E17").Value = '  -3.69%  '
Set-TextValue $ws.Range("D18") '89.19'
$ws.Range("E18").Value = '  -5.54%  '
Set-TextValue $ws.Range("D19") '0.06604'
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("E20").Value = '  -4.97%  '
$ws.Range("E21").Value = '  +0.02%  '
Set-TextValue $ws.Range("D22") '5.825'
$ws.Range("E22").Value = '  -2.46%  '
Set-TextValue $ws.Range("D23") '28.465.60'
$ws.Range("E23").Value = '  -3.85%  '
Set-TextValue $ws.Range("D24") '11.55'
$ws.Range("E24").Value = '  -3.59%  '
Set-TextValue $ws.Range("D25") '2.291'
$ws.Range("E25").Value = '  -0.49%  '
Set-TextValue $ws.Range("D26") '2.213.47'
$ws.Range("E26").Value = '  -2.44%  '
Set-TextValue $ws.Range("D27") '154.04'
$ws.Range("E27").Value = '  -3.13%  '
Set-TextValue $ws.Range("D28") '20.20'
$ws.Range("E28").Value = '  -2.56%  '
Set-TextValue $ws.Range("D29") '5.964'
$ws.Range("E29").Value = '  -6.02%  '
Set-TextValue $ws.Range("D30") '2.158'
$ws.Range("E30").Value = '  -5.99%  '
Set-TextValue $ws.Range("D31") '123.78'
$ws.Range("E31").Value = '  -3.46%  '
Set-TextValue $ws.Range("D32") '0.9888'
$ws.Range("E32").Value = '  -6.34%  '
Set-TextValue $ws.Range("D33") '0.09582'
$ws.Range("E33").Value = '  -3.82%  '
Set-TextValue $ws.Range("D34") '1.451'
$ws.Range("E34").Value = '  -7.31%  '
Set-TextValue $ws.Range("D35") '5.599'
$ws.Range("E35").Value = '  -3.98%  '
Set-TextValue $ws.Range("D36") '3.649'
$ws.Range("E36").Value = '  -4.02%  '
Set-TextValue $ws.Range("D37") '0.02331'
$ws.Range("E37").Value = '  -5.44%  '
Set-TextValue $ws.Range("D38") '8.801'
$ws.Range("E38").Value = '  -5.12%  '
Set-TextValue $ws.Range("D39") '0.06213'
$ws.Range("E39").Value = '  -3.30%  '
Set-TextValue $ws.Range("D40") '1.256'
$ws.Range("E40").Value = '  -3.79%  '
Set-TextValue $ws.Range("D41") '0.6220'
$ws.Range("E41").Value = '  -4.85%  '
Set-TextValue $ws.Range("D42") '11.14'
$ws.Range("E42").Value = '  -4.53%  '
$ws.Range("E43").Value = '  -0.19%  '
Set-TextValue $ws.Range("D44") '0.1920'
$ws.Range("E44").Value = '  -6.16%  '
Set-TextValue $ws.Range("D45") '1.340'
$ws.Range("E45").Value = '  +2.70%  '
Set-TextValue $ws.Range("D46") '0.5963'
$ws.Range("E46").Value = '  -6.12%  '
Set-TextValue $ws.Range("D47") '13.02'
$ws.Range("E47").Value = '  -3.02%  '
Set-TextValue $ws.Range("D48") '2.053'
$ws.Range("E48").Value = '  -6.48%  '
$ws.Range("E49").Value = '  -3.13%  '
Set-TextValue $ws.Range("D50") '0.00000000330'
$ws.Range("E50").Value = '  -0.95%  '
Set-TextValue $ws.Range("D51") '0.06832'
$ws.Range("E51").Value = '  -2.26%  '
